$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a "speed-up" column D = C/B (Parallel/Sequential ratio) ---
# D2 gets a plain formula; D3:D11 become a shared-formula block when
# filled down from D3, matching Excel's normal fill-down behaviour.
$ws.Range("D2").Formula = "=C2/B2"
$ws.Range("D3:D11").Formula = "=C3/B3"

# --- Fix the mislabeled "Conditions" note: it said 20 runs, should be 50 ---
$ws.Range("E4").Value = "averaged over 50 runs"

# --- View state: user zoomed in and moved the selection to E3 ---
$win = $excel.ActiveWindow
$win.Zoom = 115
$ws.Range("E3").Select()

$wb.Save()
